# ---------------------------------------------------------------------------
# Adds a new "localdb" command-category to the hidden "#system" lookup sheet
# used by the posts_and_comments worksheet's target/command drop-downs.
#
#   * inserts a brand-new column N ("localdb") on "#system", pushing every
#     existing list (macro, mail, number, pdf, rdbms, redis, sms, sound,
#     ssh, step, web, webalert, webcookie, ws, ws.async, xml) one column
#     to the right
#   * populates the new column with the localdb command set
#   * adds "localdb" to the "target" drop-down list (column A), in its
#     alphabetically-correct spot between "json" and "macro"
#   * adds two new commands (scrollElement / scrollPage) to the "web" list
#   * re-points every shifted defined name at its new range, and defines
#     the brand new "localdb" name
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# 1. Insert the new column before N (macro + everything right of it shifts
#    right by one column, taking all values/defined-name target cells with
#    it automatically).
# ---------------------------------------------------------------------------
$ws.Columns("N:N").Insert()

# ---------------------------------------------------------------------------
# 2. Populate the new "localdb" column (header + 6 commands).
# ---------------------------------------------------------------------------
$ws.Range("N1").Value() = "localdb"
$ws.Range("N2").Value() = "cloneTable(var,source,target)"
$ws.Range("N3").Value() = "dropTables(var,tables)"
$ws.Range("N4").Value() = "exportCSV(sql,output)"
$ws.Range("N5").Value() = "importRecords(var,sourceDb,sql,table)"
$ws.Range("N6").Value() = "purge(var)"
$ws.Range("N7").Value() = "runSQLs(var,sqls)"

# ---------------------------------------------------------------------------
# 3. Insert "localdb" into the "target" drop-down list (column A), between
#    "json" (row 13) and "macro" (row 14), shifting macro..xml down by one.
# ---------------------------------------------------------------------------
$ws.Range("A14").Insert(-4121)
$ws.Range("A14").Value() = "localdb"

# ---------------------------------------------------------------------------
# 4. Insert two new commands into the "web" list (column X), between
#    "saveValues(var,locator)" (row 100) and "scrollLeft(locator,pixel)"
#    (row 101).
# ---------------------------------------------------------------------------
$ws.Range("X101:X102").Insert(-4121)
$ws.Range("X101").Value() = "scrollElement(locator,xOffset,yOffset)"
$ws.Range("X102").Value() = "scrollPage(xOffset,yOffset)"

# ---------------------------------------------------------------------------
# 5. Re-point the defined names that moved, and define the new one.
# ---------------------------------------------------------------------------
function Set-NamedRange($name, $ref) {
    $n = $wb.Names.Item($name)
    $n.RefersTo() = $ref
}

Set-NamedRange "macro"     "='#system'!`$O`$2:`$O`$4"
Set-NamedRange "mail"      "='#system'!`$P`$2:`$P`$2"
Set-NamedRange "number"    "='#system'!`$Q`$2:`$Q`$16"
Set-NamedRange "pdf"       "='#system'!`$R`$2:`$R`$16"
Set-NamedRange "rdbms"     "='#system'!`$S`$2:`$S`$7"
Set-NamedRange "redis"     "='#system'!`$T`$2:`$T`$10"
Set-NamedRange "sms"       "='#system'!`$U`$2:`$U`$2"
Set-NamedRange "sound"     "='#system'!`$V`$2:`$V`$5"
Set-NamedRange "ssh"       "='#system'!`$W`$2:`$W`$9"
Set-NamedRange "step"      "='#system'!`$X`$2:`$X`$4"
Set-NamedRange "web"       "='#system'!`$Y`$2:`$Y`$127"
Set-NamedRange "webalert"  "='#system'!`$Z`$2:`$Z`$8"
Set-NamedRange "webcookie" "='#system'!`$AA`$2:`$AA`$8"
Set-NamedRange "ws"        "='#system'!`$AB`$2:`$AB`$17"
Set-NamedRange "ws.async"  "='#system'!`$AC`$2:`$AC`$8"
Set-NamedRange "xml"       "='#system'!`$AD`$2:`$AD`$21"
Set-NamedRange "target"    "='#system'!`$A`$2:`$A`$30"

$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")
